$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -19.13673333624024
$ws.Cells.Item(2, 3).Value = 1.964273483295793
$ws.Cells.Item(2, 4).Value = -19.13673333624024
$ws.Cells.Item(2, 5).Value = -19.13673333624024
$ws.Cells.Item(2, 6).Value = -19.13673333624024
$ws.Cells.Item(2, 7).Value = -19.13673333624024
$ws.Cells.Item(2, 8).Value = -19.13673333624024
$ws.Cells.Item(2, 9).Value = -19.13673333624024
$ws.Cells.Item(2, 10).Value = -19.13673333624024
$ws.Cells.Item(2, 11).Value = -19.13673333624024
$ws.Cells.Item(3, 2).Value = -19.13673333624024
$ws.Cells.Item(3, 3).Value = -19.13673333624024
$ws.Cells.Item(3, 4).Value = -19.13673333624024
$ws.Cells.Item(3, 5).Value = -19.13673333624024
$ws.Cells.Item(3, 6).Value = -19.13673333624024
$ws.Cells.Item(3, 7).Value = -19.13673333624024
$ws.Cells.Item(3, 8).Value = -19.13673333624024
$ws.Cells.Item(3, 9).Value = 1.131142618052435
$ws.Cells.Item(3, 10).Value = -19.13673333624024
$ws.Cells.Item(3, 11).Value = -19.13673333624024
$ws.Cells.Item(4, 2).Value = -19.13673333624024
$ws.Cells.Item(4, 3).Value = 2.00688218130026
$ws.Cells.Item(4, 4).Value = 1.659994701308178
$ws.Cells.Item(4, 5).Value = -19.13673333624024
$ws.Cells.Item(4, 6).Value = 3.403540645281748
$ws.Cells.Item(4, 7).Value = -19.13673333624024
$ws.Cells.Item(4, 8).Value = 1.390874434095026
$ws.Cells.Item(4, 9).Value = -19.13673333624024
$ws.Cells.Item(4, 10).Value = 1.021443638322871
$ws.Cells.Item(4, 11).Value = -19.13673333624024
$ws.Cells.Item(5, 2).Value = -19.13673333624024
$ws.Cells.Item(5, 3).Value = 1.677439940481632
$ws.Cells.Item(5, 4).Value = -19.13673333624024
$ws.Cells.Item(5, 5).Value = -19.13673333624024
$ws.Cells.Item(5, 6).Value = -19.13673333624024
$ws.Cells.Item(5, 7).Value = 2.663000843290671
$ws.Cells.Item(5, 8).Value = -19.13673333624024
$ws.Cells.Item(5, 9).Value = -19.13673333624024
$ws.Cells.Item(5, 10).Value = -19.13673333624024
$ws.Cells.Item(5, 11).Value = -19.13673333624024
$ws.Cells.Item(6, 2).Value = -19.13673333624024
$ws.Cells.Item(6, 3).Value = -19.13673333624024
$ws.Cells.Item(6, 4).Value = -19.13673333624024
$ws.Cells.Item(6, 5).Value = -19.13673333624024
$ws.Cells.Item(6, 6).Value = -19.13673333624024
$ws.Cells.Item(6, 7).Value = -19.13673333624024
$ws.Cells.Item(6, 8).Value = -19.13673333624024
$ws.Cells.Item(6, 9).Value = -19.13673333624024
$ws.Cells.Item(6, 10).Value = -19.13673333624024
$ws.Cells.Item(6, 11).Value = -19.13673333624024
$ws.Cells.Item(7, 2).Value = 2.456667361445091
$ws.Cells.Item(7, 3).Value = -19.13673333624024
$ws.Cells.Item(7, 4).Value = -19.13673333624024
$ws.Cells.Item(7, 5).Value = -19.13673333624024
$ws.Cells.Item(7, 6).Value = -19.13673333624024
$ws.Cells.Item(7, 7).Value = -19.13673333624024
$ws.Cells.Item(7, 8).Value = -19.13673333624024
$ws.Cells.Item(7, 9).Value = -19.13673333624024
$ws.Cells.Item(7, 10).Value = -19.13673333624024
$ws.Cells.Item(7, 11).Value = -19.13673333624024
$ws.Cells.Item(8, 2).Value = -19.13673333624024
$ws.Cells.Item(8, 3).Value = -19.13673333624024
$ws.Cells.Item(8, 4).Value = -19.13673333624024
$ws.Cells.Item(8, 5).Value = 1.80913410272587
$ws.Cells.Item(8, 6).Value = -19.13673333624024
$ws.Cells.Item(8, 7).Value = -19.13673333624024
$ws.Cells.Item(8, 8).Value = -19.13673333624024
$ws.Cells.Item(8, 9).Value = -19.13673333624024
$ws.Cells.Item(8, 10).Value = -19.13673333624024
$ws.Cells.Item(8, 11).Value = -19.13673333624024
$ws.Cells.Item(9, 2).Value = 3.859025017723414
$ws.Cells.Item(9, 3).Value = -19.13673333624024
$ws.Cells.Item(9, 4).Value = -19.13673333624024
$ws.Cells.Item(9, 5).Value = -19.13673333624024
$ws.Cells.Item(9, 6).Value = -19.13673333624024
$ws.Cells.Item(9, 7).Value = -19.13673333624024
$ws.Cells.Item(9, 8).Value = -19.13673333624024
$ws.Cells.Item(9, 9).Value = -19.13673333624024
$ws.Cells.Item(9, 10).Value = -19.13673333624024
$ws.Cells.Item(9, 11).Value = -19.13673333624024
$ws.Cells.Item(10, 2).Value = -19.13673333624024
$ws.Cells.Item(10, 3).Value = -19.13673333624024
$ws.Cells.Item(10, 4).Value = -19.13673333624024
$ws.Cells.Item(10, 5).Value = -19.13673333624024
$ws.Cells.Item(10, 6).Value = -19.13673333624024
$ws.Cells.Item(10, 7).Value = -19.13673333624024
$ws.Cells.Item(10, 8).Value = -19.13673333624024
$ws.Cells.Item(10, 9).Value = 1.74865993183391
$ws.Cells.Item(10, 10).Value = -19.13673333624024
$ws.Cells.Item(10, 11).Value = -19.13673333624024
$ws.Cells.Item(11, 2).Value = -19.13673333624024
$ws.Cells.Item(11, 3).Value = -19.13673333624024
$ws.Cells.Item(11, 4).Value = -19.13673333624024
$ws.Cells.Item(11, 5).Value = 2.920651950861273
$ws.Cells.Item(11, 6).Value = -19.13673333624024
$ws.Cells.Item(11, 7).Value = 2.941021535698044
$ws.Cells.Item(11, 8).Value = -19.13673333624024
$ws.Cells.Item(11, 9).Value = -19.13673333624024
$ws.Cells.Item(11, 10).Value = -19.13673333624024
$ws.Cells.Item(11, 11).Value = 4.321925717128829
$ws.Cells.Item(12, 2).Value = -19.13673333624024
$ws.Cells.Item(12, 3).Value = -19.13673333624024
$ws.Cells.Item(12, 4).Value = -19.13673333624024
$ws.Cells.Item(12, 5).Value = -19.13673333624024
$ws.Cells.Item(12, 6).Value = -19.13673333624024
$ws.Cells.Item(12, 7).Value = -19.13673333624024
$ws.Cells.Item(12, 8).Value = -19.13673333624024
$ws.Cells.Item(12, 9).Value = -19.13673333624024
$ws.Cells.Item(12, 10).Value = -19.13673333624024
$ws.Cells.Item(12, 11).Value = -19.13673333624024
$ws.Cells.Item(13, 2).Value = -19.13673333624024
$ws.Cells.Item(13, 3).Value = -19.13673333624024
$ws.Cells.Item(13, 4).Value = -19.13673333624024
$ws.Cells.Item(13, 5).Value = 2.52766924187469
$ws.Cells.Item(13, 6).Value = -19.13673333624024
$ws.Cells.Item(13, 7).Value = -19.13673333624024
$ws.Cells.Item(13, 8).Value = -19.13673333624024
$ws.Cells.Item(13, 9).Value = -19.13673333624024
$ws.Cells.Item(13, 10).Value = 1.685383873908318
$ws.Cells.Item(13, 11).Value = -19.13673333624024
$ws.Cells.Item(14, 2).Value = -19.13673333624024
$ws.Cells.Item(14, 3).Value = -19.13673333624024
$ws.Cells.Item(14, 4).Value = 1.550885826550136
$ws.Cells.Item(14, 5).Value = -19.13673333624024
$ws.Cells.Item(14, 6).Value = -19.13673333624024
$ws.Cells.Item(14, 7).Value = -19.13673333624024
$ws.Cells.Item(14, 8).Value = -19.13673333624024
$ws.Cells.Item(14, 9).Value = -19.13673333624024
$ws.Cells.Item(14, 10).Value = -19.13673333624024
$ws.Cells.Item(14, 11).Value = -19.13673333624024
$ws.Cells.Item(15, 2).Value = -19.13673333624024
$ws.Cells.Item(15, 3).Value = -19.13673333624024
$ws.Cells.Item(15, 4).Value = 1.725444661783674
$ws.Cells.Item(15, 5).Value = -19.13673333624024
$ws.Cells.Item(15, 6).Value = -19.13673333624024
$ws.Cells.Item(15, 7).Value = -19.13673333624024
$ws.Cells.Item(15, 8).Value = -19.13673333624024
$ws.Cells.Item(15, 9).Value = -19.13673333624024
$ws.Cells.Item(15, 10).Value = -19.13673333624024
$ws.Cells.Item(15, 11).Value = -19.13673333624024
$ws.Cells.Item(16, 2).Value = -19.13673333624024
$ws.Cells.Item(16, 3).Value = -19.13673333624024
$ws.Cells.Item(16, 4).Value = -19.13673333624024
$ws.Cells.Item(16, 5).Value = -19.13673333624024
$ws.Cells.Item(16, 6).Value = -19.13673333624024
$ws.Cells.Item(16, 7).Value = -19.13673333624024
$ws.Cells.Item(16, 8).Value = -19.13673333624024
$ws.Cells.Item(16, 9).Value = -19.13673333624024
$ws.Cells.Item(16, 10).Value = 1.880450921601948
$ws.Cells.Item(16, 11).Value = -19.13673333624024
$ws.Cells.Item(17, 2).Value = -19.13673333624024
$ws.Cells.Item(17, 3).Value = 2.122799039424302
$ws.Cells.Item(17, 4).Value = 1.848016790581082
$ws.Cells.Item(17, 5).Value = -19.13673333624024
$ws.Cells.Item(17, 6).Value = -19.13673333624024
$ws.Cells.Item(17, 7).Value = -19.13673333624024
$ws.Cells.Item(17, 8).Value = 2.010116839820941
$ws.Cells.Item(17, 9).Value = 2.183102483524458
$ws.Cells.Item(17, 10).Value = 2.661094792147932
$ws.Cells.Item(17, 11).Value = -19.13673333624024
$ws.Cells.Item(18, 2).Value = -19.13673333624024
$ws.Cells.Item(18, 3).Value = -19.13673333624024
$ws.Cells.Item(18, 4).Value = -19.13673333624024
$ws.Cells.Item(18, 5).Value = -19.13673333624024
$ws.Cells.Item(18, 6).Value = -19.13673333624024
$ws.Cells.Item(18, 7).Value = -19.13673333624024
$ws.Cells.Item(18, 8).Value = 1.995879978056942
$ws.Cells.Item(18, 9).Value = 1.973888242556864
$ws.Cells.Item(18, 10).Value = 2.24693562829914
$ws.Cells.Item(18, 11).Value = -19.13673333624024
$ws.Cells.Item(19, 2).Value = -19.13673333624024
$ws.Cells.Item(19, 3).Value = -19.13673333624024
$ws.Cells.Item(19, 4).Value = 2.056523096465507
$ws.Cells.Item(19, 5).Value = -19.13673333624024
$ws.Cells.Item(19, 6).Value = -19.13673333624024
$ws.Cells.Item(19, 7).Value = -19.13673333624024
$ws.Cells.Item(19, 8).Value = 1.774246754545601
$ws.Cells.Item(19, 9).Value = 1.831346303375217
$ws.Cells.Item(19, 10).Value = -19.13673333624024
$ws.Cells.Item(19, 11).Value = -19.13673333624024
$ws.Cells.Item(20, 2).Value = -19.13673333624024
$ws.Cells.Item(20, 3).Value = 1.042909742070143
$ws.Cells.Item(20, 4).Value = 1.507493517957514
$ws.Cells.Item(20, 5).Value = -19.13673333624024
$ws.Cells.Item(20, 6).Value = 3.235415683339915
$ws.Cells.Item(20, 7).Value = -19.13673333624024
$ws.Cells.Item(20, 8).Value = 1.732496757275086
$ws.Cells.Item(20, 9).Value = 1.275577037965999
$ws.Cells.Item(20, 10).Value = -19.13673333624024
$ws.Cells.Item(20, 11).Value = -19.13673333624024
$ws.Cells.Item(21, 2).Value = -19.13673333624024
$ws.Cells.Item(21, 3).Value = 1.301194454700093
$ws.Cells.Item(21, 4).Value = -19.13673333624024
$ws.Cells.Item(21, 5).Value = 1.658726914064431
$ws.Cells.Item(21, 6).Value = -19.13673333624024
$ws.Cells.Item(21, 7).Value = 2.581816680426535
$ws.Cells.Item(21, 8).Value = 1.388046520434537
$ws.Cells.Item(21, 9).Value = -19.13673333624024
$ws.Cells.Item(21, 10).Value = -19.13673333624024
$ws.Cells.Item(21, 11).Value = -19.13673333624024

Write-Output "Updated PSSM values for B2:K21"
